$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = -8.036
$ws.Range("C7").Value = -12.961
$ws.Range("E7").Value = 16.21
$ws.Range("B9").Value = 5.737
$ws.Range("E10").Value = 16.195
$ws.Range("C12").Value = -11.473
$ws.Range("E13").Value = 16.703
$ws.Range("C14").Value = -12.845
$ws.Range("D15").Value = -8.463000000000001
$ws.Range("E16").Value = 16.736
$ws.Range("B18").Value = 4.853
$ws.Range("B20").Value = 7.125
$ws.Range("E20").Value = 15.951
$ws.Range("E24").Value = 16.788
$ws.Range("C26").Value = -12.467
$ws.Range("B27").Value = 5.747999999999999
$ws.Range("C27").Value = -13.538
$ws.Range("C29").Value = -12.24
$ws.Range("D33").Value = -7.342000000000001
$ws.Range("B35").Value = 9.239000000000001
$ws.Range("D35").Value = -7.826000000000001
$ws.Range("C37").Value = -13.151
$ws.Range("C38").Value = -13.616
$ws.Range("D38").Value = -7.714
$ws.Range("E39").Value = 16.532
$ws.Range("D43").Value = -7.598999999999999
$ws.Range("D44").Value = -7.388
$ws.Range("D47").Value = -7.52
$ws.Range("E47").Value = 16.896
$ws.Range("E48").Value = 17.206
$ws.Range("C51").Value = -12.613
$ws.Range("D51").Value = -7.626
$ws.Range("C52").Value = -11.754
$ws.Range("E52").Value = 17.461
$ws.Range("C55").Value = -13.65
$ws.Range("E56").Value = 16.73
$ws.Range("D57").Value = -7.969000000000001
$ws.Range("D63").Value = -7.664999999999999
$ws.Range("B69").Value = 5.659999999999999
$ws.Range("C69").Value = -10.919
$ws.Range("C70").Value = -13.014
$ws.Range("D70").Value = -7.992999999999999
$ws.Range("B76").Value = 6.723999999999999
$ws.Range("B78").Value = 7.85
$ws.Range("C81").Value = -13.561
$ws.Range("B82").Value = 5.425
$ws.Range("B83").Value = 5.127999999999999
$ws.Range("C83").Value = -13.527
$ws.Range("E84").Value = 16.856
$ws.Range("D88").Value = -7.895
$ws.Range("B93").Value = 5.875999999999999
$ws.Range("D99").Value = -7.57
$ws.Range("E100").Value = 16.569
$ws.Range("E101").Value = 16.821
$ws.Range("C102").Value = -13.663
